# Auto-generated edit script: update cached market-price / profit values
# across multiple Leve-profit sheets (refreshed by the scheduled runner).
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3875
$ws.Range("I40").Value = 3875
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3875
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3700
$ws.Range("N40").ClearContents()
$ws.Range("H70").Value = 2190.182
$ws.Range("I70").Value = 1232.6666
$ws.Range("J70").Value = 3339.2
$ws.Range("K70").Value = 3697.9998
$ws.Range("L70").Value = 10017.6
$ws.Range("M70").Value = -3427.9998
$ws.Range("N70").Value = -10557.6
$ws.Range("H73").Value = 2190.182
$ws.Range("I73").Value = 1232.6666
$ws.Range("J73").Value = 3339.2
$ws.Range("K73").Value = 3697.9998
$ws.Range("L73").Value = 10017.6
$ws.Range("M73").Value = -2761.9998
$ws.Range("N73").Value = -11889.6
$ws.Range("H98").Value = 2468.8948
$ws.Range("I98").Value = 2161.611
$ws.Range("K98").Value = 2161.611
$ws.Range("M98").Value = -663.6109999999999
$ws.Range("H113").Value = 3991.4546
$ws.Range("I113").Value = 3066
$ws.Range("J113").Value = 5102
$ws.Range("K113").Value = 3066
$ws.Range("L113").Value = 5102
$ws.Range("M113").Value = 188
$ws.Range("N113").Value = -11610
$ws.Range("H122").Value = 2468.8948
$ws.Range("I122").Value = 2161.611
$ws.Range("K122").Value = 6484.833
$ws.Range("M122").Value = -4034.833

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1832.6666
$ws.Range("I2").Value = 1456.5714
$ws.Range("K2").Value = 1456.5714
$ws.Range("M2").Value = -1343.5714
$ws.Range("H32").Value = 29055
$ws.Range("I32").Value = 30928.893
$ws.Range("J32").Value = 17499.334
$ws.Range("K32").Value = 30928.893
$ws.Range("L32").Value = 17499.334
$ws.Range("M32").Value = -30641.893
$ws.Range("N32").Value = -18073.334
$ws.Range("H116").Value = 1832.6666
$ws.Range("I116").Value = 1456.5714
$ws.Range("K116").Value = 1456.5714
$ws.Range("M116").Value = 837.4286
$ws.Range("H122").Value = 2887.375
$ws.Range("I122").Value = 2468.3076
$ws.Range("J122").Value = 4703.3335
$ws.Range("K122").Value = 7404.9228
$ws.Range("L122").Value = 14110.0005
$ws.Range("M122").Value = -4954.9228
$ws.Range("N122").Value = -19010.0005

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1832.6666
$ws.Range("I3").Value = 1456.5714
$ws.Range("K3").Value = 1456.5714
$ws.Range("M3").Value = -1342.5714

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1393.625
$ws.Range("I6").Value = 1458.3334
$ws.Range("J6").Value = 1199.5
$ws.Range("K6").Value = 1458.3334
$ws.Range("L6").Value = 1199.5
$ws.Range("M6").Value = -1345.3334
$ws.Range("N6").Value = -1425.5
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H31").Value = 3873.875
$ws.Range("I31").Value = 2694.5
$ws.Range("K31").Value = 2694.5
$ws.Range("M31").Value = -2399.5
$ws.Range("H34").Value = 3873.875
$ws.Range("I34").Value = 2694.5
$ws.Range("K34").Value = 2694.5
$ws.Range("M34").Value = -2492.5
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H62").Value = 9999.5
$ws.Range("J62").Value = 9999.5
$ws.Range("L62").Value = 9999.5
$ws.Range("N62").Value = -11247.5
$ws.Range("H65").Value = 9999.5
$ws.Range("J65").Value = 9999.5
$ws.Range("L65").Value = 49997.5
$ws.Range("N65").Value = -56237.5
$ws.Range("H74").Value = 25000
$ws.Range("I74").Value = 20000
$ws.Range("J74").Value = 30000
$ws.Range("K74").Value = 20000
$ws.Range("L74").Value = 30000
$ws.Range("M74").Value = -19126
$ws.Range("N74").Value = -31748
$ws.Range("H77").Value = 25000
$ws.Range("I77").Value = 20000
$ws.Range("J77").Value = 30000
$ws.Range("K77").Value = 60000
$ws.Range("L77").Value = 90000
$ws.Range("M77").Value = -55632
$ws.Range("N77").Value = -98736

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 7946
$ws.Range("I139").Value = 2410
$ws.Range("J139").Value = 16250
$ws.Range("K139").Value = 7230
$ws.Range("L139").Value = 48750
$ws.Range("M139").Value = -2090
$ws.Range("N139").Value = -59030

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 15898.895
$ws.Range("J46").Value = 5505.857
$ws.Range("L46").Value = 5505.857
$ws.Range("N46").Value = -5881.857
$ws.Range("H68").Value = 4777.6665
$ws.Range("I68").Value = 10000
$ws.Range("J68").Value = 4124.875
$ws.Range("K68").Value = 10000
$ws.Range("L68").Value = 4124.875
$ws.Range("M68").Value = -9251
$ws.Range("N68").Value = -5622.875
$ws.Range("H71").Value = 4777.6665
$ws.Range("I71").Value = 10000
$ws.Range("J71").Value = 4124.875
$ws.Range("K71").Value = 50000
$ws.Range("L71").Value = 20624.375
$ws.Range("M71").Value = -46256
$ws.Range("N71").Value = -28112.375
$ws.Range("H136").Value = 9500
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 9500
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 28500
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -33600

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1231.3636
$ws.Range("I81").Value = 1244.5
$ws.Range("J81").Value = 1100
$ws.Range("K81").Value = 2489
$ws.Range("L81").Value = 2200
$ws.Range("M81").Value = -1428
$ws.Range("N81").Value = -4322
$ws.Range("H84").Value = 1231.3636
$ws.Range("I84").Value = 1244.5
$ws.Range("J84").Value = 1100
$ws.Range("K84").Value = 12445
$ws.Range("L84").Value = 11000
$ws.Range("M84").Value = -7141
$ws.Range("N84").Value = -21608
$ws.Range("H113").Value = 3107.25
$ws.Range("I113").Value = 1291.3334
$ws.Range("J113").Value = 4196.8
$ws.Range("K113").Value = 3874.0002
$ws.Range("L113").Value = 12590.4
$ws.Range("M113").Value = -1704.0002
$ws.Range("N113").Value = -16930.4
